$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, matching the style of the existing header row (A1:F1)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy formatting (font, border, alignment) from an existing header cell (F1) to G1:H1
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Add data rows for the two new columns
$ws.Range("G2").Value = 0.2668650318499809
$ws.Range("H2").Value = 0.998

$ws.Range("G3").Value = 0.2668650318499809
$ws.Range("H3").Value = 0.998
